$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.674.24'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '2.530.89'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '309.64'
$ws.Range('E5').Value = '  -1.22%  '
$ws.Range('D6').Value = '100.14'
$ws.Range('E6').Value = '  +0.84%  '
$ws.Range('E7').Value = '  -1.36%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.522'
$ws.Range('D10').Value = '35.75'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('D11').Value = '0.0803'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = '7.34'
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = '2.921.82'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').Value = '15.28'
$ws.Range('E15').Value = '  -3.59%  '
$ws.Range('D16').Value = '2.511.54'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('E17').Value = '  -3.91%  '
$ws.Range('D18').Value = '42.668.54'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.70'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('D21').Value = '12.23'
$ws.Range('E21').Value = '  -2.45%  '
$ws.Range('D22').Value = '69.29'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '242.93'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('E24').Value = '  -3.01%  '
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = '25.46'
$ws.Range('E27').Value = '  -5.65%  '
$ws.Range('D28').Value = '2.34'
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('D29').Value = '10.16'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').Value = '38.54'
$ws.Range('E30').Value = '  -3.86%  '
$ws.Range('D31').Value = '158.99'
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('D32').Value = '5.76'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('E33').Value = '  +9.78%  '
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').Value = '0.0782'
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').Value = '18.34'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').Value = '3.12'
$ws.Range('E37').Value = '  -7.15%  '
$ws.Range('D38').Value = '1.96'
$ws.Range('E38').Value = '  -7.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.110'
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('D41').Value = '4.25'
$ws.Range('E41').Value = '  +3.37%  '
$ws.Range('D42').Value = '22.53'
$ws.Range('E42').Value = '  -3.04%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '0.0299'
$ws.Range('E44').Value = '  -0.93%  '
$ws.Range('D45').Value = '3.29'
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('D46').Value = '1.994.41'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').Value = '8.89'
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('D48').Value = '2.774.89'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('D49').Value = '0.189'
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('D50').Value = '79.35'
$ws.Range('E50').Value = '  -2.93%  '
$ws.Range('D51').Value = '72.15'
$ws.Range('E51').Value = '  -3.37%  '
